$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.142.61'

# Row 3
$ws.Range("D3").Value = '1.679.32'
$ws.Range("E3").Value = '  -0.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = '''214.39'

# Row 6
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("E7").Value = '  +0.10%  '

# Row 8
$ws.Range("D8").Value = '''22.80'
$ws.Range("E8").Value = '  +6.52%  '

# Row 9
$ws.Range("D9").Value = '''0.261'
$ws.Range("E9").Value = '  +2.39%  '

# Row 10
$ws.Range("E10").Value = '  -0.36%  '

# Row 11
$ws.Range("D11").Value = '''0.0890'
$ws.Range("E11").Value = '  +0.19%  '

# Row 12
$ws.Range("D12").Value = '1.915.91'
$ws.Range("E12").Value = '  +0.01%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.685.90'
$ws.Range("E13").Value = '  -0.12%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.20'
$ws.Range("E14").Value = '  +2.14%  '

# Row 15
$ws.Range("E15").Value = '  +3.80%  '

# Row 16
$ws.Range("D16").Value = '''66.60'
$ws.Range("E16").Value = '  +0.17%  '

# Row 17
$ws.Range("D17").Value = '27.120.66'
$ws.Range("E17").Value = '  +0.38%  '

# Row 18
$ws.Range("D18").Value = '''235.05'
$ws.Range("E18").Value = '  -0.38%  '

# Row 19
$ws.Range("D19").Value = '''7.89'
$ws.Range("E19").Value = '  -3.41%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0741'
$ws.Range("E20").Value = '  +0.44%  '

# Row 21
$ws.Range("E21").Value = '  +0.09%  '

# Row 22
$ws.Range("E22").Value = '  +1.63%  '

# Row 23
$ws.Range("D23").Value = '''9.53'
$ws.Range("E23").Value = '  +2.75%  '

# Row 24
$ws.Range("D24").Value = '''2.09'
$ws.Range("E24").Value = '  -1.41%  '

# Row 25
$ws.Range("D25").Value = '''148.76'
$ws.Range("E25").Value = '  +1.55%  '

# Row 26
$ws.Range("E26").Value = '  +2.15%  '

# Row 27
$ws.Range("D27").Value = '''16.35'
$ws.Range("E27").Value = '  -0.48%  '

# Row 28
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
$ws.Range("E29").Value = '  +0.15%  '

# Row 30
$ws.Range("D30").Value = '''0.0500'
$ws.Range("E30").Value = '  +0.50%  '

# Row 31
$ws.Range("E31").Value = '  -0.47%  '

# Row 32
$ws.Range("E32").Value = '  -0.18%  '

# Row 33
$ws.Range("D33").Value = '1.539.63'
$ws.Range("E33").Value = '  +0.18%  '

# Row 34
$ws.Range("E34").Value = '  +1.13%  '

# Row 35
$ws.Range("E35").Value = '  -3.79%  '

# Row 36
$ws.Range("E36").Value = '  +2.84%  '

# Row 37
$ws.Range("D37").Value = '''0.940'
$ws.Range("E37").Value = '  +2.25%  '

# Row 38
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("E39").Value = '  -0.98%  '

# Row 40
$ws.Range("E40").Value = '  +2.15%  '

# Row 42
$ws.Range("D42").Value = '''69.44'
$ws.Range("E42").Value = '  +2.09%  '

# Row 43
$ws.Range("E43").Value = '  +0.08%  '

# Row 44
$ws.Range("E44").Value = '  -0.36%  '

# Row 45
$ws.Range("D45").Value = '1.824.16'
$ws.Range("E45").Value = '  +0.39%  '

# Row 46
$ws.Range("D46").Value = '''0.781'
$ws.Range("E46").Value = '  +0.13%  '

# Row 47
$ws.Range("D47").Value = '''89.93'
$ws.Range("E47").Value = '  -0.56%  '

# Row 48
$ws.Range("E48").Value = '  +6.05%  '

# Row 49
$ws.Range("E49").Value = '  +3.28%  '

# Row 50
$ws.Range("D50").Value = '''8.24'
$ws.Range("E50").Value = '  +3.45%  '

# Row 51
$ws.Range("E51").Value = '  -0.41%  '

